$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Features")

$ws.Range("A17").Value = 618
$ws.Range("B17").Value = 38
$ws.Range("C17").Value = 0.002629069191981492
$ws.Range("D17").Value = 0.018770602169046833
$ws.Range("E17").Value = 0.013112384595375948
$ws.Range("F17").Value = 0.000005525985514451737
$ws.Range("G17").Value = 0.00042398500789266015
$ws.Range("H17").Value = 0.0002284541730856325
